$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 4695.125
$ws.Range("I19").Value = 8761.333000000001
$ws.Range("J19").Value = 628.9167
$ws.Range("K19").Value = 8761.333000000001
$ws.Range("L19").Value = 628.9167
$ws.Range("M19").Value = -8586.333000000001
$ws.Range("N19").Value = -978.9167
$ws.Range("H76").Value = 8349.666999999999
$ws.Range("I76").Value = 11811.75
$ws.Range("J76").Value = 5580
$ws.Range("K76").Value = 11811.75
$ws.Range("L76").Value = 5580
$ws.Range("M76").Value = -11496.75
$ws.Range("N76").Value = -6210
$ws.Range("H79").Value = 8349.666999999999
$ws.Range("I79").Value = 11811.75
$ws.Range("J79").Value = 5580
$ws.Range("K79").Value = 11811.75
$ws.Range("L79").Value = 5580
$ws.Range("M79").Value = -10719.75
$ws.Range("N79").Value = -7764
$ws.Range("H137").Value = 315876.22
$ws.Range("I137").Value = 3049.923
$ws.Range("J137").Value = 824218.9399999999
$ws.Range("K137").Value = 9149.769
$ws.Range("L137").Value = 2472656.82
$ws.Range("M137").Value = -6599.769
$ws.Range("N137").Value = -2477756.82

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25190.787
$ws.Range("I32").Value = 14011.92
$ws.Range("J32").Value = 60124.75
$ws.Range("K32").Value = 14011.92
$ws.Range("L32").Value = 60124.75
$ws.Range("M32").Value = -13724.92
$ws.Range("N32").Value = -60698.75
$ws.Range("H97").Value = 2172.6
$ws.Range("I97").Value = 2234.0833
$ws.Range("K97").Value = 2234.0833
$ws.Range("M97").Value = -1738.0833
$ws.Range("H132").Value = 3450.4167
$ws.Range("I132").Value = 3294.7646
$ws.Range("J132").Value = 3828.4285
$ws.Range("K132").Value = 9884.293799999999
$ws.Range("L132").Value = 11485.2855
$ws.Range("M132").Value = -7354.293799999999
$ws.Range("N132").Value = -16545.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2340.7932
$ws.Range("I105").Value = 2164.6135
$ws.Range("K105").Value = 2164.6135
$ws.Range("M105").Value = -417.6134999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3257.9348
$ws.Range("I31").Value = 2600.4333
$ws.Range("K31").Value = 2600.4333
$ws.Range("M31").Value = -2305.4333
$ws.Range("H34").Value = 3257.9348
$ws.Range("I34").Value = 2600.4333
$ws.Range("K34").Value = 2600.4333
$ws.Range("M34").Value = -2398.4333
$ws.Range("H99").Value = 58058.555
$ws.Range("I99").Value = 73796.71000000001
$ws.Range("J99").Value = 2975
$ws.Range("K99").Value = 73796.71000000001
$ws.Range("L99").Value = 2975
$ws.Range("M99").Value = -72298.71000000001
$ws.Range("N99").Value = -5971
$ws.Range("H125").Value = 39224.75
$ws.Range("J125").Value = 39224.75
$ws.Range("L125").Value = 39224.75
$ws.Range("N125").Value = -44144.75
$ws.Range("H126").Value = 58058.555
$ws.Range("I126").Value = 73796.71000000001
$ws.Range("J126").Value = 2975
$ws.Range("K126").Value = 221390.13
$ws.Range("L126").Value = 8925
$ws.Range("M126").Value = -218920.13
$ws.Range("N126").Value = -13865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12304.889
$ws.Range("I2").Value = 18362
$ws.Range("J2").Value = 190.66667
$ws.Range("K2").Value = 110172
$ws.Range("L2").Value = 1144.00002
$ws.Range("M2").Value = -110059
$ws.Range("N2").Value = -1370.00002
$ws.Range("H3").Value = 3854.85
$ws.Range("I3").Value = 3857.3125
$ws.Range("J3").Value = 3845
$ws.Range("K3").Value = 11571.9375
$ws.Range("L3").Value = 11535
$ws.Range("M3").Value = -11459.9375
$ws.Range("N3").Value = -11759
$ws.Range("H5").Value = 735.0606
$ws.Range("I5").Value = 461.42105
$ws.Range("J5").Value = 1106.4286
$ws.Range("K5").Value = 1384.26315
$ws.Range("L5").Value = 3319.2858
$ws.Range("M5").Value = -1272.26315
$ws.Range("N5").Value = -3543.2858
$ws.Range("H14").Value = 558.75
$ws.Range("I14").Value = 558.75
$ws.Range("K14").Value = 1676.25
$ws.Range("M14").Value = -1503.25
$ws.Range("H23").Value = 446.25
$ws.Range("I23").Value = 37.5
$ws.Range("J23").Value = 477.69232
$ws.Range("K23").Value = 112.5
$ws.Range("L23").Value = 1433.07696
$ws.Range("M23").Value = 122.5
$ws.Range("N23").Value = -1903.07696
$ws.Range("H68").Value = 662.5
$ws.Range("J68").Value = 685.7143
$ws.Range("L68").Value = 2057.1429
$ws.Range("N68").Value = -3679.1429
$ws.Range("H71").Value = 662.5
$ws.Range("J71").Value = 685.7143
$ws.Range("L71").Value = 6171.428699999999
$ws.Range("N71").Value = -14283.4287
$ws.Range("H98").Value = 18433.334
$ws.Range("I98").Value = 300
$ws.Range("J98").Value = 27500
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 82500
$ws.Range("M98").Value = 598
$ws.Range("N98").Value = -85496
$ws.Range("H113").Value = 478.36365
$ws.Range("I113").Value = 491.63635
$ws.Range("J113").Value = 438.54544
$ws.Range("K113").Value = 1474.90905
$ws.Range("L113").Value = 1315.63632
$ws.Range("M113").Value = 695.09095
$ws.Range("N113").Value = -5655.63632
$ws.Range("H122").Value = 627.6818
$ws.Range("I122").Value = 429.25
$ws.Range("J122").Value = 865.8
$ws.Range("K122").Value = 3863.25
$ws.Range("L122").Value = 7792.2
$ws.Range("M122").Value = -1413.25
$ws.Range("N122").Value = -12692.2
$ws.Range("H132").Value = 6313.0713
$ws.Range("I132").Value = 2643.4546
$ws.Range("J132").Value = 19768.334
$ws.Range("K132").Value = 23791.0914
$ws.Range("L132").Value = 177915.006
$ws.Range("M132").Value = -21261.0914
$ws.Range("N132").Value = -182975.006
$ws.Range("H135").Value = 735.0606
$ws.Range("I135").Value = 461.42105
$ws.Range("J135").Value = 1106.4286
$ws.Range("K135").Value = 4152.78945
$ws.Range("L135").Value = 9957.857399999999
$ws.Range("M135").Value = -1617.78945
$ws.Range("N135").Value = -15027.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5031.4287
$ws.Range("I80").Value = 6111.4287
$ws.Range("J80").Value = 2871.4285
$ws.Range("K80").Value = 6111.4287
$ws.Range("L80").Value = 2871.4285
$ws.Range("M80").Value = -5113.4287
$ws.Range("N80").Value = -4867.4285
$ws.Range("H83").Value = 5031.4287
$ws.Range("I83").Value = 6111.4287
$ws.Range("J83").Value = 2871.4285
$ws.Range("K83").Value = 30557.1435
$ws.Range("L83").Value = 14357.1425
$ws.Range("M83").Value = -25565.1435
$ws.Range("N83").Value = -24341.1425
$ws.Range("H107").Value = 7916.154
$ws.Range("I107").Value = 7916.154
$ws.Range("K107").Value = 7916.154
$ws.Range("M107").Value = -5996.154
$ws.Range("H122").Value = 1351.8966
$ws.Range("I122").Value = 1234.6522
$ws.Range("J122").Value = 1801.3334
$ws.Range("K122").Value = 3703.9566
$ws.Range("L122").Value = 5404.0002
$ws.Range("M122").Value = -1253.9566
$ws.Range("N122").Value = -10304.0002
$ws.Range("H132").Value = 6335.5
$ws.Range("I132").Value = 10201.714
$ws.Range("J132").Value = 4253.6924
$ws.Range("K132").Value = 30605.142
$ws.Range("L132").Value = 12761.0772
$ws.Range("M132").Value = -28075.142
$ws.Range("N132").Value = -17821.0772

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1728.85
$ws.Range("I7").Value = 1641.8889
$ws.Range("J7").Value = 1800
$ws.Range("K7").Value = 1641.8889
$ws.Range("L7").Value = 1800
$ws.Range("M7").Value = -1529.8889
$ws.Range("N7").Value = -2024
$ws.Range("H40").Value = 2336.0715
$ws.Range("I40").Value = 2528.5715
$ws.Range("J40").Value = 2143.5715
$ws.Range("K40").Value = 2528.5715
$ws.Range("L40").Value = 2143.5715
$ws.Range("M40").Value = -2392.5715
$ws.Range("N40").Value = -2415.5715
$ws.Range("H61").Value = 1677.6086
$ws.Range("I61").Value = 1649.6666
$ws.Range("J61").Value = 1778.2
$ws.Range("K61").Value = 1649.6666
$ws.Range("L61").Value = 1778.2
$ws.Range("M61").Value = -1447.6666
$ws.Range("N61").Value = -2182.2
$ws.Range("H82").Value = 2610
$ws.Range("I82").Value = 2377.1428
$ws.Range("J82").Value = 3153.3333
$ws.Range("K82").Value = 2377.1428
$ws.Range("L82").Value = 3153.3333
$ws.Range("M82").Value = -2016.1428
$ws.Range("N82").Value = -3875.3333
$ws.Range("H85").Value = 2610
$ws.Range("I85").Value = 2377.1428
$ws.Range("J85").Value = 3153.3333
$ws.Range("K85").Value = 2377.1428
$ws.Range("L85").Value = 3153.3333
$ws.Range("M85").Value = -1129.1428
$ws.Range("N85").Value = -5649.3333
$ws.Range("H113").Value = 1677.6086
$ws.Range("I113").Value = 1649.6666
$ws.Range("J113").Value = 1778.2
$ws.Range("K113").Value = 1649.6666
$ws.Range("L113").Value = 1778.2
$ws.Range("M113").Value = 520.3334
$ws.Range("N113").Value = -6118.2
$ws.Range("H126").Value = 1728.85
$ws.Range("I126").Value = 1641.8889
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4925.6667
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2455.6667
$ws.Range("N126").Value = -10340
$ws.Range("H132").Value = 9813657
$ws.Range("I132").Value = 3684.4
$ws.Range("J132").Value = 23827902
$ws.Range("K132").Value = 11053.2
$ws.Range("L132").Value = 71483706
$ws.Range("M132").Value = -8523.200000000001
$ws.Range("N132").Value = -71488766

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 65278.938
$ws.Range("I122").Value = 113995.89
$ws.Range("J122").Value = 2642.8572
$ws.Range("K122").Value = 341987.67
$ws.Range("L122").Value = 7928.571599999999
$ws.Range("M122").Value = -339537.67
$ws.Range("N122").Value = -12828.5716
$ws.Range("H126").Value = 892
$ws.Range("I126").Value = 742.37036
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 2227.11108
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = 242.8889199999999
$ws.Range("N126").Value = -10040
$ws.Range("H132").Value = 2036.6923
$ws.Range("I132").Value = 1114.1666
$ws.Range("J132").Value = 4112.375
$ws.Range("K132").Value = 3342.4998
$ws.Range("L132").Value = 12337.125
$ws.Range("M132").Value = -812.4998000000001
$ws.Range("N132").Value = -17397.125
